$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 3 data (A3, B3) first, keeping the cells present with default style
$ws.Range("A3").Value = ""
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("B3").Value = ""
$ws.Range("B3").VerticalAlignment = -4108

# Update B2 text from 배고파 to YAMANG93
$ws.Range("B2").Value = "YAMANG93"

# Add new header cell C1 = VISIBLE
$ws.Range("C1").Value = "VISIBLE"

# Move selection to C1
$ws.Range("C1").Select()
